# Applies: insert new sheet "2022-Q1" (fund-holding detail, after "2021-Q2"
# and before "总计"), and update the "总计" (totals) sheet with a new
# first data row summarising the 2022-Q1 quarter, pushing the existing
# 2021-Q2 total row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q2".
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$q1Sheet = $wb.Worksheets.Add($null, $firstSheet)
$q1Sheet.Name = "2022-Q1"

# Headers (row 1, columns B:H - column A is left blank, matching the
# layout used by the other holdings sheets in this workbook).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q1Sheet.Cells.Item(1, 2 + $i)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: index, fund code, fund name, fund size, total stock
# position, position ratio, held market value (billion yuan), position
# rank.
$rows = @(
    @(0, "519029", "华夏稳增混合", "10.92", "92.99", "5.66", "0.6181", 3),
    @(1, "000480", "东方红新动力灵活配置混合", "15.38", "72.90", "2.81", "0.4322", 10),
    @(2, "001564", "东方红京东大数据灵活配置混合", "11.40", "69.58", "2.20", "0.2508", 10)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q1Sheet.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    # Text-valued columns B:G - force Text format so values such as
    # "10.92" / "0.6181" are preserved exactly instead of becoming
    # numbers.
    for ($c = 2; $c -le 7; $c++) {
        $cell = $q1Sheet.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
    }

    # Rank column (H) is numeric.
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new summary row for
#    2022-Q1 above the existing 2021-Q2 row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$a2 = $totalSheet.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 1.3

# The row that was pushed down to row 3 (previously row 2, "2021-Q2")
# keeps its data, but its index value needs to move from 0 to 1.
$totalSheet.Cells.Item(3, 1).Value = 1

# Restore the originally active sheet/tab.
$firstSheet.Activate()
